$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 51, shifting existing rows 51..102 down to 52..103.
$ws.Rows("51:51").Insert()

# Populate the newly inserted row 51 with the new data record.
$ws.Range("A51").Value = 1
$ws.Range("B51").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C51").Value = "Arica y Parinacota"
$ws.Range("D51").Value = 44977
$ws.Range("E51").Value = 15
$ws.Range("F51").Value = "Fruta"
$ws.Range("G51").Value = 100103
$ws.Range("H51").Value = "Frutos de hueso (carozo)"
$ws.Range("I51").Value = 100103006
$ws.Range("J51").Value = "Nectarín"
$ws.Range("K51").Value = "Artic Sprite"
$ws.Range("L51").Value = "Primera"
$ws.Range("M51").Value = 550
$ws.Range("N51").Value = 22000
$ws.Range("O51").Value = 23000
$ws.Range("P51").Value = 22636
$ws.Range("Q51").Value = '$/caja 20 kilos granel'
$ws.Range("R51").Value = "Región de O'Higgins"
$ws.Range("S51").Value = 1132
$ws.Range("T51").Value = 20

# Match the date-number-format style (s="2") used by column D elsewhere in the sheet.
$ws.Range("D51").NumberFormat = "YYYY-MM-DD HH:MM:SS"
